$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 20:40"

# 2. Swap country names that changed rank order (their stat rows keep the
#    row position of their rank, so swap which country name sits in each
#    row, then below we set the new stats for each row).

# Marruecos overtakes Uzbekistan (rows 64/65)
$ws.Range("A64").Value = "Marruecos"
$ws.Range("A65").Value = "Uzbekistan"

# Sudan overtakes Bulgaria (rows 81/82)
$ws.Range("A81").Value = "Sudan"
$ws.Range("A82").Value = "Bulgaria"

# Grecia overtakes Hungria (rows 105/106)
$ws.Range("A105").Value = "Grecia"
$ws.Range("A106").Value = "Hungria"

# 3. Update numeric stats (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected rows.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4736014
$ws.Range("C4").Value = 30125
$ws.Range("D4").Value = 2338613
$ws.Range("E4").Value = 2240123
$ws.Range("G4").Value = 531
$ws.Range("H4").Value = 157278

# Row 6 - India
$ws.Range("B6").Value = 1751836
$ws.Range("C6").Value = 54782
$ws.Range("D6").Value = 1146828
$ws.Range("E6").Value = 567605
$ws.Range("G6").Value = 852
$ws.Range("H6").Value = 37403

# Row 20 - Turquia
$ws.Range("B20").Value = 231869
$ws.Range("C20").Value = 996
$ws.Range("D20").Value = 215516
$ws.Range("E20").Value = 10643
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 5710

# Row 64 - now Marruecos (new stats)
$ws.Range("B64").Value = 25015
$ws.Range("C64").Value = 693
$ws.Range("D64").Value = 17960
$ws.Range("E64").Value = 6688
$ws.Range("G64").Value = 14
$ws.Range("H64").Value = 367

# Row 65 - now Uzbekistan (keeps its previous stats, unchanged values)
$ws.Range("B65").Value = 24569
$ws.Range("C65").Value = 560
$ws.Range("D65").Value = 14916
$ws.Range("E65").Value = 9508
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 145

# Row 81 - now Sudan (new stats)
$ws.Range("B81").Value = 11738
$ws.Range("C81").Value = 94
$ws.Range("D81").Value = 6137
$ws.Range("E81").Value = 4849
$ws.Range("G81").Value = 6
$ws.Range("H81").Value = 752

# Row 82 - now Bulgaria (keeps its previous stats, unchanged values)
$ws.Range("B82").Value = 11690
$ws.Range("D82").Value = 6319
$ws.Range("E82").Value = 4988
$ws.Range("H82").Value = 383

# Row 89 - Guayana Francesa
$ws.Range("B89").Value = 7857
$ws.Range("C89").Value = 58
$ws.Range("D89").Value = 6531
$ws.Range("E89").Value = 1283

# Row 105 - now Grecia (new stats)
$ws.Range("B105").Value = 4587
$ws.Range("C105").Value = 110
$ws.Range("D105").Value = 1374
$ws.Range("E105").Value = 3007
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 206

# Row 106 - now Hungria (keeps its previous stats, unchanged values)
$ws.Range("B106").Value = 4526
$ws.Range("C106").Value = 21
$ws.Range("D106").Value = 3364
$ws.Range("E106").Value = 565
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 597

# Row 135 - Yemen
$ws.Range("B135").Value = 1730
$ws.Range("C135").Value = 2
$ws.Range("E135").Value = 374
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 494
